$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of data produced by the solution (graphing of the produced solution)
# Columns: A=Node number, B=X coord, C=Y coord, D=Demand, E=Type
$data = @(
    @(4, -1, 2, 10, "L"),
    @(5, -4, -1, 20, "L"),
    @(6, 0, -2, 10, "B"),
    @(7, -2, -1, 30, "L"),
    @(8, -4, -2, 10, "B")
)

$startRow = 6
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}

$ws.Range("E11").Select()
